{"js": "// The document contains a single table of arithmetic problems (20 rows x 5\n// columns = 100 cells). The edit replaces the text of every cell with a new\n// expression, in row-major (reading) order. The pairs below are\n// [oldText, newText] taken straight from the authoritative diff, in the\n// exact order the cells appear in the document.\nconst REPLACEMENTS = [[\"41-26=15\",\"69+28=97\"],[\"24+48=72\",\"48+44=92\"],[\"60-45=15\",\"39+22=61\"],[\"48+49=97\",\"23-16=7\"],[\"92-25=67\",\"49+24=73\"],[\"44-37=7\",\"39+22=61\"],[\"5+17=22\",\"17+65=82\"],[\"39+9=48\",\"8+88=96\"],[\"81-36=45\",\"43-35=8\"],[\"5+38=43\",\"95-39=56\"],[\"49+16=65\",\"24+57=81\"],[\"92-57=35\",\"6+18=24\"],[\"70-49=21\",\"6+28=34\"],[\"81-72=9\",\"92-6=86\"],[\"19+18=37\",\"92-89=3\"],[\"94-58=36\",\"74-46=28\"],[\"70-41=29\",\"62-16=46\"],[\"24+47=71\",\"61-14=47\"],[\"80-78=2\",\"25+49=74\"],[\"23-4=19\",\"8+64=72\"],[\"78+9=87\",\"56+9=65\"],[\"14+77=91\",\"7+37=44\"],[\"46+15=61\",\"43-34=9\"],[\"18+55=73\",\"22-6=16\"],[\"50-17=33\",\"73-38=35\"],[\"71-53=18\",\"19+5=24\"],[\"56-9=47\",\"45+46=91\"],[\"17+17=34\",\"41-27=14\"],[\"41-12=29\",\"30-4=26\"],[\"64-25=39\",\"32-7=25\"],[\"14+69=83\",\"60-18=42\"],[\"46-28=18\",\"70-6=64\"],[\"69+5=74\",\"97-89=8\"],[\"38+28=66\",\"70-12=58\"],[\"35-26=9\",\"73-55=18\"],[\"18+23=41\",\"17+18=35\"],[\"29+14=43\",\"40-31=9\"],[\"80-74=6\",\"54-36=18\"],[\"17+35=52\",\"90-31=59\"],[\"91-67=24\",\"20-7=13\"],[\"50-29=21\",\"6+77=83\"],[\"37+47=84\",\"72-39=33\"],[\"29+66=95\",\"65-47=18\"],[\"62-44=18\",\"56+6=62\"],[\"96-48=48\",\"12+29=41\"],[\"69+24=93\",\"90-87=3\"],[\"27+65=92\",\"23-17=6\"],[\"74+7=81\",\"48+9=57\"],[\"92-16=76\",\"49+43=92\"],[\"17+5=22\",\"91-26=65\"],[\"18+17=35\",\"80-23=57\"],[\"26-8=18\",\"47-38=9\"],[\"15+66=81\",\"51-35=16\"],[\"83-9=74\",\"28+19=47\"],[\"90-87=3\",\"92-15=77\"],[\"95-76=19\",\"60-1=59\"],[\"85-77=8\",\"64+7=71\"],[\"49+18=67\",\"4+79=83\"],[\"7+5=12\",\"77-28=49\"],[\"96-88=8\",\"94-85=9\"],[\"95-66=29\",\"22+69=91\"],[\"29+38=67\",\"70-17=53\"],[\"7+87=94\",\"94-16=78\"],[\"69+4=73\",\"3+49=52\"],[\"45-6=39\",\"47+18=65\"],[\"59+2=61\",\"48+14=62\"],[\"18+23=41\",\"76+8=84\"],[\"27-19=8\",\"21-5=16\"],[\"42-27=15\",\"73-65=8\"],[\"57+4=61\",\"9+9=18\"],[\"32-5=27\",\"39+53=92\"],[\"15+49=64\",\"95-9=86\"],[\"82-15=67\",\"20-19=1\"],[\"3+18=21\",\"14+9=23\"],[\"83-57=26\",\"92-63=29\"],[\"70-14=56\",\"24+27=51\"],[\"13+19=32\",\"28+54=82\"],[\"56+7=63\",\"52-26=26\"],[\"32-15=17\",\"6+16=22\"],[\"17+35=52\",\"48+28=76\"],[\"70-19=51\",\"7+9=16\"],[\"13+18=31\",\"40-31=9\"],[\"20-18=2\",\"90-27=63\"],[\"84-57=27\",\"57+15=72\"],[\"81-15=66\",\"74+8=82\"],[\"17+45=62\",\"61-39=22\"],[\"91-33=58\",\"71-28=43\"],[\"26+35=61\",\"85-29=56\"],[\"55+37=92\",\"89+3=92\"],[\"73-25=48\",\"70-36=34\"],[\"68+5=73\",\"84-47=37\"],[\"95-78=17\",\"90-8=82\"],[\"93-69=24\",\"45-28=17\"],[\"63-5=58\",\"53-48=5\"],[\"47+4=51\",\"52-39=13\"],[\"69+18=87\",\"72-35=37\"],[\"79+9=88\",\"52-13=39\"],[\"49+5=54\",\"86-8=78\"],[\"16+16=32\",\"91-78=13\"],[\"29+54=83\",\"65-46=19\"]];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount, values\");\nawait context.sync();\n\nconst rowCount = table.rowCount;\nconst colCount = table.values[0].length;\n\n// Build the full replacement matrix, walking the REPLACEMENTS list in the\n// same row-major order the cells occupy. Flatten the existing values first\n// so we can sanity-check the \"old\" side of each pair against what's\n// actually in the document before swapping in the \"new\" side.\nconst flatCurrent = [];\nfor (let r = 0; r < rowCount; r++) {\n  for (let c = 0; c < colCount; c++) {\n    flatCurrent.push(table.values[r][c]);\n  }\n}\n\nconst flatNew = flatCurrent.slice();\nfor (let i = 0; i < REPLACEMENTS.length && i < flatNew.length; i++) {\n  const [oldText, newText] = REPLACEMENTS[i];\n  // Only swap in the new text when the cell still holds the expected old\n  // text; otherwise leave the cell untouched rather than corrupt data.\n  if (flatCurrent[i] === oldText) {\n    flatNew[i] = newText;\n  }\n}\n\nconst newMatrix = [];\nfor (let r = 0; r < rowCount; r++) {\n  newMatrix.push(flatNew.slice(r * colCount, r * colCount + colCount));\n}\n\n// Single bulk write of the whole table's text values.\ntable.values = newMatrix;\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n$rows = $t.Rows.Count\n$cols = $t.Columns.Count\n\n# Expected current text and replacement text for every cell, in row-major\n# (reading) order -- taken directly from the authoritative diff.\n$oldValues = @(\n    \"41-26=15\",\n    \"24+48=72\",\n    \"60-45=15\",\n    \"48+49=97\",\n    \"92-25=67\",\n    \"44-37=7\",\n    \"5+17=22\",\n    \"39+9=48\",\n    \"81-36=45\",\n    \"5+38=43\",\n    \"49+16=65\",\n    \"92-57=35\",\n    \"70-49=21\",\n    \"81-72=9\",\n    \"19+18=37\",\n    \"94-58=36\",\n    \"70-41=29\",\n    \"24+47=71\",\n    \"80-78=2\",\n    \"23-4=19\",\n    \"78+9=87\",\n    \"14+77=91\",\n    \"46+15=61\",\n    \"18+55=73\",\n    \"50-17=33\",\n    \"71-53=18\",\n    \"56-9=47\",\n    \"17+17=34\",\n    \"41-12=29\",\n    \"64-25=39\",\n    \"14+69=83\",\n    \"46-28=18\",\n    \"69+5=74\",\n    \"38+28=66\",\n    \"35-26=9\",\n    \"18+23=41\",\n    \"29+14=43\",\n    \"80-74=6\",\n    \"17+35=52\",\n    \"91-67=24\",\n    \"50-29=21\",\n    \"37+47=84\",\n    \"29+66=95\",\n    \"62-44=18\",\n    \"96-48=48\",\n    \"69+24=93\",\n    \"27+65=92\",\n    \"74+7=81\",\n    \"92-16=76\",\n    \"17+5=22\",\n    \"18+17=35\",\n    \"26-8=18\",\n    \"15+66=81\",\n    \"83-9=74\",\n    \"90-87=3\",\n    \"95-76=19\",\n    \"85-77=8\",\n    \"49+18=67\",\n    \"7+5=12\",\n    \"96-88=8\",\n    \"95-66=29\",\n    \"29+38=67\",\n    \"7+87=94\",\n    \"69+4=73\",\n    \"45-6=39\",\n    \"59+2=61\",\n    \"18+23=41\",\n    \"27-19=8\",\n    \"42-27=15\",\n    \"57+4=61\",\n    \"32-5=27\",\n    \"15+49=64\",\n    \"82-15=67\",\n    \"3+18=21\",\n    \"83-57=26\",\n    \"70-14=56\",\n    \"13+19=32\",\n    \"56+7=63\",\n    \"32-15=17\",\n    \"17+35=52\",\n    \"70-19=51\",\n    \"13+18=31\",\n    \"20-18=2\",\n    \"84-57=27\",\n    \"81-15=66\",\n    \"17+45=62\",\n    \"91-33=58\",\n    \"26+35=61\",\n    \"55+37=92\",\n    \"73-25=48\",\n    \"68+5=73\",\n    \"95-78=17\",\n    \"93-69=24\",\n    \"63-5=58\",\n    \"47+4=51\",\n    \"69+18=87\",\n    \"79+9=88\",\n    \"49+5=54\",\n    \"16+16=32\",\n    \"29+54=83\"\n)\n$newValues = @(\n    \"69+28=97\",\n    \"48+44=92\",\n    \"39+22=61\",\n    \"23-16=7\",\n    \"49+24=73\",\n    \"39+22=61\",\n    \"17+65=82\",\n    \"8+88=96\",\n    \"43-35=8\",\n    \"95-39=56\",\n    \"24+57=81\",\n    \"6+18=24\",\n    \"6+28=34\",\n    \"92-6=86\",\n    \"92-89=3\",\n    \"74-46=28\",\n    \"62-16=46\",\n    \"61-14=47\",\n    \"25+49=74\",\n    \"8+64=72\",\n    \"56+9=65\",\n    \"7+37=44\",\n    \"43-34=9\",\n    \"22-6=16\",\n    \"73-38=35\",\n    \"19+5=24\",\n    \"45+46=91\",\n    \"41-27=14\",\n    \"30-4=26\",\n    \"32-7=25\",\n    \"60-18=42\",\n    \"70-6=64\",\n    \"97-89=8\",\n    \"70-12=58\",\n    \"73-55=18\",\n    \"17+18=35\",\n    \"40-31=9\",\n    \"54-36=18\",\n    \"90-31=59\",\n    \"20-7=13\",\n    \"6+77=83\",\n    \"72-39=33\",\n    \"65-47=18\",\n    \"56+6=62\",\n    \"12+29=41\",\n    \"90-87=3\",\n    \"23-17=6\",\n    \"48+9=57\",\n    \"49+43=92\",\n    \"91-26=65\",\n    \"80-23=57\",\n    \"47-38=9\",\n    \"51-35=16\",\n    \"28+19=47\",\n    \"92-15=77\",\n    \"60-1=59\",\n    \"64+7=71\",\n    \"4+79=83\",\n    \"77-28=49\",\n    \"94-85=9\",\n    \"22+69=91\",\n    \"70-17=53\",\n    \"94-16=78\",\n    \"3+49=52\",\n    \"47+18=65\",\n    \"48+14=62\",\n    \"76+8=84\",\n    \"21-5=16\",\n    \"73-65=8\",\n    \"9+9=18\",\n    \"39+53=92\",\n    \"95-9=86\",\n    \"20-19=1\",\n    \"14+9=23\",\n    \"92-63=29\",\n    \"24+27=51\",\n    \"28+54=82\",\n    \"52-26=26\",\n    \"6+16=22\",\n    \"48+28=76\",\n    \"7+9=16\",\n    \"40-31=9\",\n    \"90-27=63\",\n    \"57+15=72\",\n    \"74+8=82\",\n    \"61-39=22\",\n    \"71-28=43\",\n    \"85-29=56\",\n    \"89+3=92\",\n    \"70-36=34\",\n    \"84-47=37\",\n    \"90-8=82\",\n    \"45-28=17\",\n    \"53-48=5\",\n    \"52-39=13\",\n    \"72-35=37\",\n    \"52-13=39\",\n    \"86-8=78\",\n    \"91-78=13\",\n    \"65-46=19\"\n)\n\n$idx = 0\nfor ($r = 1; $r -le $rows; $r++) {\n  for ($c = 1; $c -le $cols; $c++) {\n    if ($idx -lt $oldValues.Length) {\n      $cell = $t.Cell($r, $c)\n      # A whole-cell Range.Text includes the trailing end-of-cell marker\n      # (CR + BEL), so strip it before comparing/assigning plain text.\n      $current = $cell.Range.Text.TrimEnd([char]13, [char]7)\n      if ($current -eq $oldValues[$idx]) {\n        $cell.Range.Text = $newValues[$idx]\n      }\n    }\n    $idx++\n  }\n}\n"}
